# Update the "Förändrad" date column (C) for rows 2-8 from 45212 to 45221
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
